$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.889.92"
$ws.Range("E2").Value = "  +2.21%  "

$ws.Range("D3").Value = "3.061.69"
$ws.Range("E3").Value = "  +2.66%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.56"
$ws.Range("E5").Value = "  +5.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.30"
$ws.Range("E6").Value = "  +6.09%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  +5.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.64"
$ws.Range("E9").Value = "  +5.78%  "

$ws.Range("E10").Value = "  +7.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.370"
$ws.Range("E11").Value = "  +5.84%  "

$ws.Range("E12").Value = "  +2.02%  "

$ws.Range("D13").Value = "3.589.10"
$ws.Range("E13").Value = "  +2.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.33"
$ws.Range("E14").Value = "  +8.20%  "

$ws.Range("E15").Value = "  +16.29%  "

$ws.Range("D16").Value = "57.877.44"
$ws.Range("E16").Value = "  +2.19%  "

$ws.Range("E17").Value = "  +8.59%  "

$ws.Range("D18").Value = "3.062.88"
$ws.Range("E18").Value = "  +2.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.18"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.22"
$ws.Range("E20").Value = "  +5.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "341.42"
$ws.Range("E21").Value = "  +5.19%  "

$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("E23").Value = "  +7.32%  "

$ws.Range("E24").Value = "  +5.32%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "0.0₃0978"
$ws.Range("E25").Value = "  +9.54%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("E26").Value = "  +5.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("E28").Value = "  +7.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.41"
$ws.Range("E29").Value = "  +9.74%  "

$ws.Range("E30").Value = "  +6.93%  "

$ws.Range("E31").Value = "  +7.24%  "

$ws.Range("E32").Value = "  +4.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.79"
$ws.Range("E33").Value = "  +8.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "157.14"
$ws.Range("E34").Value = "  +2.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.98"
$ws.Range("E35").Value = "  +7.03%  "

$ws.Range("E36").Value = "  +4.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.30"
$ws.Range("E37").Value = "  +12.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0706"
$ws.Range("E38").Value = "  +5.43%  "

$ws.Range("D39").Value = "3.098.01"
$ws.Range("E39").Value = "  +2.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.90"
$ws.Range("E40").Value = "  +3.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.94"
$ws.Range("E41").Value = "  +10.95%  "

$ws.Range("E42").Value = "  -0.16%  "

$ws.Range("E43").Value = "  +4.31%  "

$ws.Range("E44").Value = "  +5.86%  "

$ws.Range("D45").Value = "2.341.89"
$ws.Range("E45").Value = "  +5.38%  "

$ws.Range("E46").Value = "  +3.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.01"
$ws.Range("E47").Value = "  +2.98%  "

$ws.Range("E48").Value = "  +5.63%  "

$ws.Range("E49").Value = "  +4.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.20"
$ws.Range("E50").Value = "  +6.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0901"
$ws.Range("E51").Value = "  +6.26%  "
